# Generate Report for Handoff
# Updates the zh-cn and de-de handoff-status sheets: flips the status from
# "Handoff transform failed" to "Ready for handoff", records the freshly
# produced xlf handoff file (with a hyperlink to it) and its handoff
# timestamp, and switches the per-row "Handoff Reason" from "Ignored" to
# "Include" now that the file is ready to go out. The Overview rollup is
# kept in sync with the same status.

$wb = $excel.ActiveWorkbook

$commitBase = "https://github.com/OpenLocalizationTest/oltest/blob/2110df0eb61378830358a339a1b35907049fc90d"

function Set-HandoffReady {
    param([string]$SheetName, [string]$XlfFileName, [string]$HandoffDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # Status: no longer failed -- ready to hand off.
    $ws.Range("B2").Value = "Ready for handoff"

    # Latest Handoff File (column C): new xlf produced by the transform,
    # linked back to the file in the source repo.
    $ws.Hyperlinks.Add($ws.Range("C2"), "$commitBase/e2e/$XlfFileName", [Type]::Missing, [Type]::Missing, $XlfFileName) | Out-Null
    $ws.Range("C2").Font.Underline = $true
    $ws.Range("C2").Font.Color = 15570276

    # Latest Handoff Datetime (column D): stamp of this handoff run.
    $ws.Range("D2").Value = $HandoffDateTime

    # Handoff Reason (column H): file is included in this handoff batch.
    $ws.Range("H2").Value = "Include"
}

Set-HandoffReady "zh-cn" "88e15722-ee89-4119-8452-37afc5243a5f.a1530fc6d372cad3cc89744f5077efe799337559.zh-cn.xlf" "2016-01-27 02:58:41"
Set-HandoffReady "de-de" "88e15722-ee89-4119-8452-37afc5243a5f.a1530fc6d372cad3cc89744f5077efe799337559.de-de.xlf" "2016-01-27 02:58:53"

# The Overview sheet rolls up the per-language status for each source file
# (column B = zh-cn, column C = de-de) -- keep it in sync with the same
# "Ready for handoff" status just recorded on the language sheets.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
